{"js": "// Replace the date line and each two-digit multiplication answer in the\n// table with the updated values from the commit. Every old text value in\n// this document is unique, so a case-sensitive exact search on the old\n// text unambiguously identifies the run to update.\nconst replacements = [\n  [\"2025-09-23 Tuesday\", \"2025-09-24 Wednesday\"],\n  [\"83\u00d731=2573\", \"92\u00d754=4968\"],\n  [\"18\u00d772=1296\", \"26\u00d754=1404\"],\n  [\"70\u00d796=6720\", \"50\u00d766=3300\"],\n  [\"48\u00d768=3264\", \"47\u00d712=564\"],\n  [\"30\u00d762=1860\", \"52\u00d760=3120\"],\n  [\"47\u00d765=3055\", \"13\u00d774=962\"],\n  [\"97\u00d748=4656\", \"95\u00d745=4275\"],\n  [\"28\u00d755=1540\", \"82\u00d748=3936\"],\n  [\"52\u00d773=3796\", \"69\u00d761=4209\"],\n  [\"89\u00d724=2136\", \"43\u00d719=817\"],\n  [\"15\u00d752=780\", \"66\u00d749=3234\"],\n  [\"52\u00d713=676\", \"93\u00d759=5487\"],\n  [\"81\u00d731=2511\", \"62\u00d726=1612\"],\n  [\"86\u00d760=5160\", \"57\u00d725=1425\"],\n  [\"13\u00d719=247\", \"31\u00d785=2635\"],\n  [\"54\u00d757=3078\", \"99\u00d735=3465\"],\n  [\"41\u00d737=1517\", \"31\u00d785=2635\"],\n  [\"79\u00d788=6952\", \"90\u00d769=6210\"],\n  [\"68\u00d757=3876\", \"86\u00d755=4730\"],\n  [\"97\u00d734=3298\", \"26\u00d775=1950\"],\n  [\"60\u00d766=3960\", \"99\u00d765=6435\"],\n  [\"41\u00d771=2911\", \"61\u00d766=4026\"],\n  [\"17\u00d757=969\", \"65\u00d737=2405\"],\n  [\"42\u00d787=3654\", \"93\u00d712=1116\"],\n  [\"60\u00d744=2640\", \"74\u00d736=2664\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the date line and every two-digit x two-digit multiplication\n# answer in the table to the values from the new day's worksheet.\n# Every \"old\" text value in this document is unique, so a case-sensitive\n# Find/Replace on the exact old text unambiguously targets the right run.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-09-23 Tuesday\", \"2025-09-24 Wednesday\"),\n    @(\"83\u00d731=2573\", \"92\u00d754=4968\"),\n    @(\"18\u00d772=1296\", \"26\u00d754=1404\"),\n    @(\"70\u00d796=6720\", \"50\u00d766=3300\"),\n    @(\"48\u00d768=3264\", \"47\u00d712=564\"),\n    @(\"30\u00d762=1860\", \"52\u00d760=3120\"),\n    @(\"47\u00d765=3055\", \"13\u00d774=962\"),\n    @(\"97\u00d748=4656\", \"95\u00d745=4275\"),\n    @(\"28\u00d755=1540\", \"82\u00d748=3936\"),\n    @(\"52\u00d773=3796\", \"69\u00d761=4209\"),\n    @(\"89\u00d724=2136\", \"43\u00d719=817\"),\n    @(\"15\u00d752=780\", \"66\u00d749=3234\"),\n    @(\"52\u00d713=676\", \"93\u00d759=5487\"),\n    @(\"81\u00d731=2511\", \"62\u00d726=1612\"),\n    @(\"86\u00d760=5160\", \"57\u00d725=1425\"),\n    @(\"13\u00d719=247\", \"31\u00d785=2635\"),\n    @(\"54\u00d757=3078\", \"99\u00d735=3465\"),\n    @(\"41\u00d737=1517\", \"31\u00d785=2635\"),\n    @(\"79\u00d788=6952\", \"90\u00d769=6210\"),\n    @(\"68\u00d757=3876\", \"86\u00d755=4730\"),\n    @(\"97\u00d734=3298\", \"26\u00d775=1950\"),\n    @(\"60\u00d766=3960\", \"99\u00d765=6435\"),\n    @(\"41\u00d771=2911\", \"61\u00d766=4026\"),\n    @(\"17\u00d757=969\", \"65\u00d737=2405\"),\n    @(\"42\u00d787=3654\", \"93\u00d712=1116\"),\n    @(\"60\u00d744=2640\", \"74\u00d736=2664\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
